$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.281.43"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "3.746.55"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D7").Value = "3.744.88"
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "4.383.84"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "3.763.28"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").Value = "69.307.25"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +17.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000154"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").Value = "3.897.61"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "3.683.54"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "436.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.801.00"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0354"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.87%  "
